$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.96"

$ws.Range("D3").Value = "'23.03"

$ws.Range("D4").Value = "'5.418"

$ws.Range("D5").Value = "'0.05991"

$ws.Range("D6").Value = "'3.394"

$ws.Range("D8").Value = "'0.9263"

$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01123"
$ws.Range("E9").Value = "8OneONEBestin24h"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1433"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07434"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03378"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03039"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09359"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.934"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001597"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04813"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("D18").Value = "'0.005534"

$ws.Range("D19").Value = "'0.004157"

$ws.Range("D20").Value = "'0.0009836"

$ws.Range("D23").Value = "'6.457"

$ws.Range("D24").Value = "'2.185"

$ws.Range("D40").Value = "'0.03952"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006152"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.006839"

$ws.Range("D45").Value = "'0.00005204"

$ws.Range("D48").Value = "'1.090"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

